# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# (crypto price/volume refresh + dogwifhat/PolygonEcosystemToken row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Price" column writes to remain plain text (matches source data,
# which stores prices such as "576.49" or "70.522.61" as text, not numbers).
$priceCells = @("D2","D3","D5","D6","D8","D9","D10","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D27","D29","D30","D31","D32","D33","D35","D36","D37","D39","D40","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "70.522.61"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.552.12"
$ws.Range("E3").Value = "  -4.87%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "576.49"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").Value = "170.12"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.510"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").Value = "2.553.47"
$ws.Range("E9").Value = "  -4.87%  "
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").Value = "4.82"
$ws.Range("E13").Value = "  -3.38%  "
$ws.Range("D14").Value = "3.021.15"
$ws.Range("E14").Value = "  -4.91%  "
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "70.427.52"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "25.12"
$ws.Range("E17").Value = "  -4.12%  "
$ws.Range("D18").Value = "2.553.97"
$ws.Range("E18").Value = "  -4.87%  "
$ws.Range("D19").Value = "11.68"
$ws.Range("E19").Value = "  -4.57%  "
$ws.Range("D20").Value = "7.68"
$ws.Range("E20").Value = "  -5.98%  "
$ws.Range("D21").Value = "361.58"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("D22").Value = "3.95"
$ws.Range("E22").Value = "  -5.16%  "
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "69.92"
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("E26").Value = "  -5.46%  "
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  -4.55%  "
$ws.Range("E28").Value = "  -4.75%  "
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").Value = "486.49"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "0.116"
$ws.Range("E36").Value = "  +6.61%  "
$ws.Range("D37").Value = "156.97"
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").Value = "18.82"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").Value = "1.32"
$ws.Range("E40").Value = "  -4.21%  "
$ws.Range("D43").Value = "4.77"
$ws.Range("E43").Value = "  -4.70%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "0.321"
$ws.Range("E44").Value = "  -3.41%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("D46").Value = "38.44"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("D47").Value = "145.57"
$ws.Range("E47").Value = "  -6.92%  "
$ws.Range("D48").Value = "3.56"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").Value = "0.531"
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("D50").Value = "1.63"
$ws.Range("E50").Value = "  -6.49%  "
$ws.Range("E51").Value = "  -1.90%  "
